$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 25000736
$ws.Range("I19").Value = 637.5714
$ws.Range("J19").Value = 38462330
$ws.Range("K19").Value = 637.5714
$ws.Range("L19").Value = 38462330
$ws.Range("M19").Value = -462.5714
$ws.Range("N19").Value = -38462680
$ws.Range("H28").Value = 736.8182
$ws.Range("I28").Value = 588.125
$ws.Range("J28").Value = 1133.3334
$ws.Range("K28").Value = 588.125
$ws.Range("L28").Value = 1133.3334
$ws.Range("M28").Value = -103.125
$ws.Range("N28").Value = -2103.3334
$ws.Range("H111").Value = 2024.1538
$ws.Range("I111").Value = 1457.1428
$ws.Range("J111").Value = 2685.6667
$ws.Range("K111").Value = 4371.428400000001
$ws.Range("L111").Value = 8057.000100000001
$ws.Range("M111").Value = -1304.428400000001
$ws.Range("N111").Value = -14191.0001
$ws.Range("H116").Value = 3862.5
$ws.Range("I116").Value = 3835
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 3835
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -393
$ws.Range("N116").Value = -10884
$ws.Range("H129").Value = 1165.4744
$ws.Range("I129").Value = 649
$ws.Range("J129").Value = 1298.758
$ws.Range("K129").Value = 1947
$ws.Range("L129").Value = 3896.274
$ws.Range("M129").Value = 3053
$ws.Range("N129").Value = -13896.274
$ws.Range("H137").Value = 1587.8695
$ws.Range("I137").Value = 1901.9259
$ws.Range("J137").Value = 1141.579
$ws.Range("K137").Value = 5705.7777
$ws.Range("L137").Value = 3424.737
$ws.Range("M137").Value = -3155.7777
$ws.Range("N137").Value = -8524.737000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1028.0526
$ws.Range("I2").Value = 602.2
$ws.Range("J2").Value = 2625
$ws.Range("K2").Value = 602.2
$ws.Range("L2").Value = 2625
$ws.Range("M2").Value = -489.2
$ws.Range("N2").Value = -2851
$ws.Range("H32").Value = 8418.591
$ws.Range("I32").Value = 7374.1484
$ws.Range("K32").Value = 7374.1484
$ws.Range("M32").Value = -7087.1484
$ws.Range("H116").Value = 1028.0526
$ws.Range("I116").Value = 602.2
$ws.Range("J116").Value = 2625
$ws.Range("K116").Value = 602.2
$ws.Range("L116").Value = 2625
$ws.Range("M116").Value = 1691.8
$ws.Range("N116").Value = -7213
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1028.0526
$ws.Range("I3").Value = 602.2
$ws.Range("J3").Value = 2625
$ws.Range("K3").Value = 602.2
$ws.Range("L3").Value = 2625
$ws.Range("M3").Value = -488.2
$ws.Range("N3").Value = -2853
$ws.Range("H86").Value = 1106.762
$ws.Range("I86").Value = 978.36365
$ws.Range("J86").Value = 1248
$ws.Range("K86").Value = 978.36365
$ws.Range("L86").Value = 1248
$ws.Range("M86").Value = 144.63635
$ws.Range("N86").Value = -3494
$ws.Range("H89").Value = 1106.762
$ws.Range("I89").Value = 978.36365
$ws.Range("J89").Value = 1248
$ws.Range("K89").Value = 4891.81825
$ws.Range("L89").Value = 6240
$ws.Range("M89").Value = 724.1817499999997
$ws.Range("N89").Value = -17472
$ws.Range("H94").Value = 800.25
$ws.Range("I94").Value = 776.2273
$ws.Range("J94").Value = 888.3333
$ws.Range("K94").Value = 776.2273
$ws.Range("L94").Value = 888.3333
$ws.Range("M94").Value = -325.2273
$ws.Range("N94").Value = -1790.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8416.083000000001
$ws.Range("I7").Value = 12547.625
$ws.Range("J7").Value = 153
$ws.Range("K7").Value = 12547.625
$ws.Range("L7").Value = 153
$ws.Range("M7").Value = -12434.625
$ws.Range("N7").Value = -379
$ws.Range("H16").Value = 6777.75
$ws.Range("I16").Value = 7703.6665
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 7703.6665
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -7416.6665
$ws.Range("N16").Value = -4574
$ws.Range("H62").Value = 9472.5
$ws.Range("I62").Value = 2467.9167
$ws.Range("J62").Value = 51500
$ws.Range("K62").Value = 2467.9167
$ws.Range("L62").Value = 51500
$ws.Range("M62").Value = -1843.9167
$ws.Range("N62").Value = -52748
$ws.Range("H65").Value = 9472.5
$ws.Range("I65").Value = 2467.9167
$ws.Range("J65").Value = 51500
$ws.Range("K65").Value = 12339.5835
$ws.Range("L65").Value = 257500
$ws.Range("M65").Value = -9219.583500000001
$ws.Range("N65").Value = -263740
$ws.Range("H113").Value = 6777.75
$ws.Range("I113").Value = 7703.6665
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 7703.6665
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -5533.6665
$ws.Range("N113").Value = -8340
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 119.666664
$ws.Range("I14").Value = 119.666664
$ws.Range("K14").Value = 358.999992
$ws.Range("M14").Value = -185.999992
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3761877.5
$ws.Range("I102").Value = 6495136.5
$ws.Range("J102").Value = 3646.75
$ws.Range("K102").Value = 6495136.5
$ws.Range("L102").Value = 3646.75
$ws.Range("M102").Value = -6493514.5
$ws.Range("N102").Value = -6890.75
$ws.Range("H132").Value = 1738716.5
$ws.Range("I132").Value = 2316592.8
$ws.Range("J132").Value = 5087.6665
$ws.Range("K132").Value = 6949778.399999999
$ws.Range("L132").Value = 15262.9995
$ws.Range("M132").Value = -6947248.399999999
$ws.Range("N132").Value = -20322.9995
$ws.Range("H133").Value = 50780
$ws.Range("J133").Value = 50780
$ws.Range("L133").Value = 50780
$ws.Range("N133").Value = -60900
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4763.4287
$ws.Range("I61").Value = 4908.8
$ws.Range("J61").Value = 4400
$ws.Range("K61").Value = 4908.8
$ws.Range("L61").Value = 4400
$ws.Range("M61").Value = -4706.8
$ws.Range("N61").Value = -4804
$ws.Range("H68").Value = 1642.8125
$ws.Range("I68").Value = 1444.5834
$ws.Range("J68").Value = 2237.5
$ws.Range("K68").Value = 1444.5834
$ws.Range("L68").Value = 2237.5
$ws.Range("M68").Value = -695.5834
$ws.Range("N68").Value = -3735.5
$ws.Range("H71").Value = 1642.8125
$ws.Range("I71").Value = 1444.5834
$ws.Range("J71").Value = 2237.5
$ws.Range("K71").Value = 7222.916999999999
$ws.Range("L71").Value = 11187.5
$ws.Range("M71").Value = -3478.916999999999
$ws.Range("N71").Value = -18675.5
$ws.Range("H113").Value = 4763.4287
$ws.Range("I113").Value = 4908.8
$ws.Range("J113").Value = 4400
$ws.Range("K113").Value = 4908.8
$ws.Range("L113").Value = 4400
$ws.Range("M113").Value = -2738.8
$ws.Range("N113").Value = -8740
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 729.6667
$ws.Range("I107").Value = 667.4286
$ws.Range("J107").Value = 854.1429000000001
$ws.Range("K107").Value = 2002.2858
$ws.Range("L107").Value = 2562.4287
$ws.Range("M107").Value = -82.28579999999988
$ws.Range("N107").Value = -6402.4287
